# Applies the cryptos price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.392.17'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.848.29'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D5").Value = '''240.62'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '''0.6330'
$ws.Range("E6").Value = '  -3.38%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.07593'
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("D9").Value = '''0.2971'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = '''24.52'
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").Value = '2.592.00'
$ws.Range("E11").Value = '  +40.27%  '
$ws.Range("D12").Value = '''0.07723'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '''0.6855'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '''82.92'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Value = '''0.000009944'
$ws.Range("E16").Value = '  +3.99%  '
$ws.Range("D17").Value = '''6.185'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '29.423.08'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '''231.30'
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D22").Value = '''7.601'
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '''154.76'
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").Value = '''0.1398'
$ws.Range("D26").Value = '''8.446'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Value = '''17.67'
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '''1.472'
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").Value = '''0.05817'
$ws.Range("E29").Value = '  -3.47%  '
$ws.Range("D30").Value = '''1.267'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("D31").Value = '''4.124'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '''1.866'
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("D35").Value = '''0.7174'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = '2.709.81'
$ws.Range("E36").Value = '  +34.60%  '
$ws.Range("D37").Value = '''2.598'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '1.249.11'
$ws.Range("E38").Value = '  +3.91%  '
$ws.Range("D39").Value = '''2.795'
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").Value = '''0.01811'
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").Value = '''0.9050'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("E42").Value = '  -2.76%  '
$ws.Range("D43").Value = '''0.9995'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '''101.41'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = '''67.23'
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").Value = '''7.319'
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("D47").Value = '''9.164'
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = '''0.4012'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").Value = '''1.692'
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("E51").Value = '  +0.07%  '
